$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Remove the "integer / what_is_your_age?" row (row 3) entirely.
$ws.Range("A3:F3").EntireRow.Delete()

# Update row 2 from the "text / what_is_your_name?" question to an "image" question.
$ws.Range("A2").Value = "image"
$ws.Range("B2").Value = "this_is_image"
$ws.Range("C2").Value = "This is image"

# Add the new trailing columns (hint, default, guidance_hint, hxl) to the header row.
$ws.Range("G1").Value = "hint"
$ws.Range("H1").Value = "default"
$ws.Range("I1").Value = "guidance_hint"
$ws.Range("J1").Value = "hxl"

# Keep the new columns present (empty) on row 2 as well, matching the header's extent.
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
